$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 14801
$ws.Range("I20").Value = 14801
$ws.Range("K20").Value = 14801
$ws.Range("M20").Value = -14571
$ws.Range("H35").Value = 14801
$ws.Range("I35").Value = 14801
$ws.Range("K35").Value = 14801
$ws.Range("M35").Value = -14422
$ws.Range("H44").Value = 14681.667
$ws.Range("J44").Value = 19500
$ws.Range("L44").Value = 19500
$ws.Range("N44").Value = -20424
$ws.Range("H103").Value = 500.2857
$ws.Range("I103").Value = 443.42856
$ws.Range("J103").Value = 557.1429000000001
$ws.Range("K103").Value = 1330.28568
$ws.Range("L103").Value = 1671.4287
$ws.Range("M103").Value = -744.28568
$ws.Range("N103").Value = -2843.4287
$ws.Range("H116").Value = 8698261
$ws.Range("I116").Value = 18183810
$ws.Range("J116").Value = 3174.9167
$ws.Range("K116").Value = 18183810
$ws.Range("L116").Value = 3174.9167
$ws.Range("M116").Value = -18180368
$ws.Range("N116").Value = -10058.9167
$ws.Range("H125").Value = 500499.5
$ws.Range("I125").Value = 500499.5
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 4504495.5
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -4502035.5
$ws.Range("N125").ClearContents()
$ws.Range("H131").Value = 8599.317999999999
$ws.Range("I131").Value = 2060
$ws.Range("J131").Value = 9631.842000000001
$ws.Range("K131").Value = 6180
$ws.Range("L131").Value = 28895.526
$ws.Range("M131").Value = -1140
$ws.Range("N131").Value = -38975.526
$ws.Range("H135").Value = 593.1429000000001
$ws.Range("I135").Value = 562.8
$ws.Range("J135").Value = 1200
$ws.Range("K135").Value = 5065.2
$ws.Range("L135").Value = 10800
$ws.Range("M135").Value = -2530.2
$ws.Range("N135").Value = -15870
$ws.Range("H137").Value = 1357.5
$ws.Range("I137").Value = 943.75
$ws.Range("J137").Value = 2350.5
$ws.Range("K137").Value = 2831.25
$ws.Range("L137").Value = 7051.5
$ws.Range("M137").Value = -281.25
$ws.Range("N137").Value = -12151.5
$ws.Range("H141").Value = 6257.3447
$ws.Range("I141").Value = 2922.0833
$ws.Range("J141").Value = 22266.6
$ws.Range("K141").Value = 8766.249899999999
$ws.Range("L141").Value = 66799.79999999999
$ws.Range("M141").Value = -3586.249899999999
$ws.Range("N141").Value = -77159.79999999999

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H57").Value = 4249.5
$ws.Range("I57").Value = 4249.5
$ws.Range("K57").Value = 4249.5
$ws.Range("M57").Value = -3765.5
$ws.Range("H88").Value = 1913.8572
$ws.Range("I88").Value = 1400
$ws.Range("J88").Value = 2119.4
$ws.Range("K88").Value = 1400
$ws.Range("L88").Value = 2119.4
$ws.Range("M88").Value = -994
$ws.Range("N88").Value = -2931.4
$ws.Range("H91").Value = 1913.8572
$ws.Range("I91").Value = 1400
$ws.Range("J91").Value = 2119.4
$ws.Range("K91").Value = 1400
$ws.Range("L91").Value = 2119.4
$ws.Range("M91").Value = 4
$ws.Range("N91").Value = -4927.4
$ws.Range("H123").Value = 24295.834
$ws.Range("J123").Value = 24295.834
$ws.Range("L123").Value = 24295.834
$ws.Range("N123").Value = -34095.834
$ws.Range("H132").Value = 2194.08
$ws.Range("I132").Value = 1773.6765
$ws.Range("J132").Value = 3087.4375
$ws.Range("K132").Value = 5321.029500000001
$ws.Range("L132").Value = 9262.3125
$ws.Range("M132").Value = -2791.029500000001
$ws.Range("N132").Value = -14322.3125

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H109").Value = 22556.5
$ws.Range("J109").Value = 22556.5
$ws.Range("L109").Value = 22556.5
$ws.Range("N109").Value = -25330.5
$ws.Range("H140").Value = 190000
$ws.Range("J140").Value = 190000
$ws.Range("L140").Value = 190000
$ws.Range("N140").Value = -200360

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 45000
$ws.Range("I54").Value = 45000
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 45000
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -44342
$ws.Range("N54").ClearContents()
$ws.Range("H94").Value = 4474.75
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 4474.75
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 4474.75
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -5376.75
$ws.Range("H134").Value = 1157.0807
$ws.Range("I134").Value = 939.3333
$ws.Range("J134").Value = 1903.6428
$ws.Range("K134").Value = 2817.9999
$ws.Range("L134").Value = 5710.928400000001
$ws.Range("M134").Value = -282.9998999999998
$ws.Range("N134").Value = -10780.9284

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 15502.857
$ws.Range("J82").Value = 15502.857
$ws.Range("L82").Value = 46508.571
$ws.Range("N82").Value = -47320.571
$ws.Range("H85").Value = 15502.857
$ws.Range("J85").Value = 15502.857
$ws.Range("L85").Value = 46508.571
$ws.Range("N85").Value = -49316.571
$ws.Range("H132").Value = 2494
$ws.Range("J132").Value = 2792
$ws.Range("L132").Value = 25128
$ws.Range("N132").Value = -30188

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 20000000
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H102").Value = 3177.9333
$ws.Range("I102").Value = 3090.5454
$ws.Range("J102").Value = 3228.5264
$ws.Range("K102").Value = 3090.5454
$ws.Range("L102").Value = 3228.5264
$ws.Range("M102").Value = -1468.5454
$ws.Range("N102").Value = -6472.526400000001
$ws.Range("H109").Value = 9279.799999999999
$ws.Range("J109").Value = 9279.799999999999
$ws.Range("L109").Value = 9279.799999999999
$ws.Range("N109").Value = -11359.8
$ws.Range("H122").Value = 2856.5881
$ws.Range("I122").Value = 2089.4443
$ws.Range("J122").Value = 3719.625
$ws.Range("K122").Value = 6268.3329
$ws.Range("L122").Value = 11158.875
$ws.Range("M122").Value = -3818.3329
$ws.Range("N122").Value = -16058.875
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3342.8667
$ws.Range("I7").Value = 3230.7273
$ws.Range("J7").Value = 3651.25
$ws.Range("K7").Value = 3230.7273
$ws.Range("L7").Value = 3651.25
$ws.Range("M7").Value = -3118.7273
$ws.Range("N7").Value = -3875.25
$ws.Range("H16").Value = 1527.6
$ws.Range("I16").Value = 1334.4445
$ws.Range("J16").Value = 1817.3334
$ws.Range("K16").Value = 1334.4445
$ws.Range("L16").Value = 1817.3334
$ws.Range("M16").Value = -1164.4445
$ws.Range("N16").Value = -2157.3334
$ws.Range("H40").Value = 4398.933
$ws.Range("I40").Value = 5110.5
$ws.Range("K40").Value = 5110.5
$ws.Range("M40").Value = -4974.5
$ws.Range("H122").Value = 19571026
$ws.Range("I122").Value = 31256350
$ws.Range("J122").Value = 13338853
$ws.Range("K122").Value = 93769050
$ws.Range("L122").Value = 40016559
$ws.Range("M122").Value = -93766600
$ws.Range("N122").Value = -40021459
$ws.Range("H126").Value = 3342.8667
$ws.Range("I126").Value = 3230.7273
$ws.Range("J126").Value = 3651.25
$ws.Range("K126").Value = 9692.1819
$ws.Range("L126").Value = 10953.75
$ws.Range("M126").Value = -7222.1819
$ws.Range("N126").Value = -15893.75
$ws.Range("H132").Value = 3704.0417
$ws.Range("I132").Value = 4709
$ws.Range("J132").Value = 2699.0833
$ws.Range("K132").Value = 14127
$ws.Range("L132").Value = 8097.249899999999
$ws.Range("M132").Value = -11597
$ws.Range("N132").Value = -13157.2499
$ws.Range("H136").Value = 1920.5646
$ws.Range("I136").Value = 1463.585
$ws.Range("J136").Value = 4611.6665
$ws.Range("K136").Value = 4390.755
$ws.Range("L136").Value = 13834.9995
$ws.Range("M136").Value = -1840.755
$ws.Range("N136").Value = -18934.9995

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4342.857
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 4480
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 4480
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -5728
$ws.Range("H65").Value = 4342.857
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 4480
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 22400
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -28640
$ws.Range("H113").Value = 577.17645
$ws.Range("I113").Value = 434.13333
$ws.Range("J113").Value = 1650
$ws.Range("K113").Value = 1302.39999
$ws.Range("L113").Value = 4950
$ws.Range("M113").Value = 867.6000100000001
$ws.Range("N113").Value = -9290
$ws.Range("H123").Value = 28860.5
$ws.Range("J123").Value = 28860.5
$ws.Range("L123").Value = 28860.5
$ws.Range("N123").Value = -38660.5
$ws.Range("H132").Value = 887.6129
$ws.Range("I132").Value = 626.86
$ws.Range("J132").Value = 1974.0834
$ws.Range("K132").Value = 1880.58
$ws.Range("L132").Value = 5922.2502
$ws.Range("M132").Value = 649.4200000000001
$ws.Range("N132").Value = -10982.2502
$ws.Range("H135").Value = 49000
$ws.Range("J135").Value = 49000
$ws.Range("L135").Value = 49000
$ws.Range("N135").Value = -59140
$ws.Range("H136").Value = 851.80853
$ws.Range("I136").Value = 881.2222
$ws.Range("J136").Value = 190
$ws.Range("K136").Value = 2643.6666
$ws.Range("L136").Value = 570
$ws.Range("M136").Value = -93.66660000000002
$ws.Range("N136").Value = -5670

Write-Host "Applied all Asura_Profits updates"